# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: Francia ---
$ws.Range("B7").Value = 162100
$ws.Range("C7").Value = 612
$ws.Range("E7").Value = 94341

# --- Row 14: Brasil ---
$ws.Range("B14").Value = 60311
$ws.Range("C14").Value = 1115
$ws.Range("E14").Value = 27034
$ws.Range("G14").Value = 72
$ws.Range("H14").Value = 4117

# --- Row 16: Canada ---
$ws.Range("B16").Value = 45800
$ws.Range("C16").Value = 446
$ws.Range("E16").Value = 26393
$ws.Range("G16").Value = 26
$ws.Range("H16").Value = 2491

# --- Rows 31 & 32: swap Japon / Pakistan (Pakistan now listed first with updated stats,
#     Japon moves down keeping its previous stats) ---
$ws.Range("A31").Value = "Pakistan"
$ws.Range("B31").Value = 13328
$ws.Range("C31").Value = 605
$ws.Range("D31").Value = 2936
$ws.Range("E31").Value = 10111
$ws.Range("F31").Value = 111
$ws.Range("G31").Value = 12
$ws.Range("H31").Value = 281

$ws.Range("A32").Value = "Japon"
$ws.Range("B32").Value = 13231
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 1656
$ws.Range("E32").Value = 11215
$ws.Range("F32").Value = 287
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 360

# --- Row 44: Noruega ---
$ws.Range("B44").Value = 7511
$ws.Range("C44").Value = 18
$ws.Range("E44").Value = 7278
$ws.Range("F44").Value = 49

# --- Row 107: Georgia ---
$ws.Range("D107").Value = 149
$ws.Range("E107").Value = 330
